# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.953.57'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '2.734.90'
$ws.Range('E3').Value = '  +3.45%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '607.08'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '170.27'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +7.05%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.551'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.71%  '
$ws.Range('D9').Value = '2.733.37'
$ws.Range('E9').Value = '  +3.44%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.149'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.49%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.370'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +5.45%  '
$ws.Range('E12').Value = '  +1.51%  '
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('E14').Value = '  +3.62%  '
$ws.Range('D15').Value = '3.232.94'
$ws.Range('E15').Value = '  +3.42%  '
$ws.Range('E16').Value = '  +2.75%  '
$ws.Range('D17').Value = '68.846.26'
$ws.Range('E17').Value = '  +1.13%  '
$ws.Range('D18').Value = '2.691.01'
$ws.Range('E18').Value = '  +1.90%  '
$ws.Range('E19').Value = '  +5.17%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '376.92'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +5.21%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.74'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.96%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.55'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.25%  '
$ws.Range('E23').Value = '  +5.95%  '
$ws.Range('E24').Value = '  +3.32%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '73.79'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.12%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.23'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.33%  '
$ws.Range('D28').Value = '2.875.43'
$ws.Range('E28').Value = '  +3.45%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0000106'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +3.28%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '588.02'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +5.43%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.42'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.48%  '
$ws.Range('E33').Value = '  +4.79%  '
$ws.Range('E34').Value = '  +6.31%  '
$ws.Range('E35').Value = '  +4.41%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '162.30'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.83%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '20.04'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.66%  '
$ws.Range('E40').Value = '  +3.78%  '
$ws.Range('E41').Value = '  +3.65%  '
$ws.Range('E42').Value = '  +3.80%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.68'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.05%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '17.99'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.04%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '41.17'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.82%  '
$ws.Range('D47').Value = '0.0₆0311'
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.608'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +6.57%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '156.29'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.97'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +4.75%  '
$ws.Range('B51').Value = 'Optimism'
$ws.Range('C51').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.81'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +7.45%  '
